$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'24.368.75"
$ws.Cells.Item(2, 5).Value = "  +9.03%  "

$ws.Cells.Item(3, 4).Value = "'1.676.81"
$ws.Cells.Item(3, 5).Value = "  +4.62%  "

$ws.Cells.Item(4, 4).Value = "'1.003"
$ws.Cells.Item(4, 5).Value = "  -0.39%  "

$ws.Cells.Item(5, 4).Value = "'306.60"
$ws.Cells.Item(5, 5).Value = "  +6.13%  "

$ws.Cells.Item(6, 4).Value = "'0.9976"
$ws.Cells.Item(6, 5).Value = "  +0.06%  "

$ws.Cells.Item(7, 4).Value = "'0.3706"
$ws.Cells.Item(7, 5).Value = "  -0.19%  "

$ws.Cells.Item(8, 4).Value = "'0.3433"
$ws.Cells.Item(8, 5).Value = "  +1.19%  "

$ws.Cells.Item(9, 4).Value = "'48.20"
$ws.Cells.Item(9, 5).Value = "  +13.82%  "

$ws.Cells.Item(10, 4).Value = "'1.177"
$ws.Cells.Item(10, 5).Value = "  +2.91%  "

$ws.Cells.Item(11, 4).Value = "'0.07245"
$ws.Cells.Item(11, 5).Value = "  +2.64%  "

$ws.Cells.Item(12, 4).Value = "'0.9999"
$ws.Cells.Item(12, 5).Value = "  -0.13%  "

$ws.Cells.Item(13, 4).Value = "'20.33"
$ws.Cells.Item(13, 5).Value = "  +2.42%  "

$ws.Cells.Item(14, 4).Value = "'6.092"
$ws.Cells.Item(14, 5).Value = "  +2.51%  "

$ws.Cells.Item(15, 4).Value = "'6.733"
$ws.Cells.Item(15, 5).Value = "  +1.01%  "

$ws.Cells.Item(16, 4).Value = "'1.679.80"
$ws.Cells.Item(16, 5).Value = "  +5.19%  "

$ws.Cells.Item(17, 4).Value = "'0.00001107"
$ws.Cells.Item(17, 5).Value = "  +2.01%  "

$ws.Cells.Item(18, 4).Value = "'0.9981"
$ws.Cells.Item(18, 5).Value = "  +0.21%  "

$ws.Cells.Item(19, 4).Value = "'0.06720"
$ws.Cells.Item(19, 5).Value = "  +1.32%  "

$ws.Cells.Item(20, 4).Value = "'81.01"
$ws.Cells.Item(20, 5).Value = "  +3.42%  "

$ws.Cells.Item(21, 4).Value = "'16.41"
$ws.Cells.Item(21, 5).Value = "  +1.27%  "

$ws.Cells.Item(22, 4).Value = "'6.085"
$ws.Cells.Item(22, 5).Value = "  +0.48%  "

$ws.Cells.Item(23, 4).Value = "'11.93"
$ws.Cells.Item(23, 5).Value = "  +1.04%  "

$ws.Cells.Item(24, 4).Value = "'24.317.72"
$ws.Cells.Item(24, 5).Value = "  +8.98%  "

$ws.Cells.Item(25, 4).Value = "'2.428"
$ws.Cells.Item(25, 5).Value = "  +1.29%  "

$ws.Cells.Item(26, 4).Value = "'3.364"
$ws.Cells.Item(26, 5).Value = "  -12.12%  "

$ws.Cells.Item(27, 4).Value = "'2.657"
$ws.Cells.Item(27, 5).Value = "  +6.02%  "

$ws.Cells.Item(28, 4).Value = "'152.30"
$ws.Cells.Item(28, 5).Value = "  +0.68%  "

$ws.Cells.Item(29, 4).Value = "'19.54"
$ws.Cells.Item(29, 5).Value = "  -0.83%  "

$ws.Cells.Item(30, 4).Value = "'1.862.69"
$ws.Cells.Item(30, 5).Value = "  +4.85%  "

$ws.Cells.Item(31, 4).Value = "'126.79"
$ws.Cells.Item(31, 5).Value = "  +4.74%  "

$ws.Cells.Item(32, 4).Value = "'6.292"
$ws.Cells.Item(32, 5).Value = "  +4.29%  "

$ws.Cells.Item(33, 4).Value = "'4.032"
$ws.Cells.Item(33, 5).Value = "  -3.37%  "

$ws.Cells.Item(34, 4).Value = "'0.9666"
$ws.Cells.Item(34, 5).Value = "  +1.82%  "

$ws.Cells.Item(35, 4).Value = "'1.737"
$ws.Cells.Item(35, 5).Value = "  +7.92%  "

$ws.Cells.Item(36, 4).Value = "'0.08456"
$ws.Cells.Item(36, 5).Value = "  +2.35%  "

$ws.Cells.Item(37, 4).Value = "'8.976"
$ws.Cells.Item(37, 5).Value = "  +3.57%  "

$ws.Cells.Item(38, 4).Value = "'12.27"
$ws.Cells.Item(38, 5).Value = "  +3.79%  "

$ws.Cells.Item(39, 4).Value = "'0.06429"
$ws.Cells.Item(39, 5).Value = "  +4.40%  "

$ws.Cells.Item(40, 4).Value = "'5.327"
$ws.Cells.Item(40, 5).Value = "  +0.36%  "

$ws.Cells.Item(41, 4).Value = "'0.02332"
$ws.Cells.Item(41, 5).Value = "  +5.07%  "

$ws.Cells.Item(42, 4).Value = "'1.262"
$ws.Cells.Item(42, 5).Value = "  +1.65%  "

$ws.Cells.Item(43, 4).Value = "'0.2106"
$ws.Cells.Item(43, 5).Value = "  +3.87%  "

$ws.Cells.Item(44, 4).Value = "'0.6152"
$ws.Cells.Item(44, 5).Value = "  +3.44%  "

$ws.Cells.Item(45, 5).Value = "  +0.13%  "

$ws.Cells.Item(46, 4).Value = "'3.774"
$ws.Cells.Item(46, 5).Value = "  +2.94%  "

$ws.Cells.Item(47, 4).Value = "'13.02"
$ws.Cells.Item(47, 5).Value = "  -0.97%  "

$ws.Cells.Item(48, 4).Value = "'0.5931"
$ws.Cells.Item(48, 5).Value = "  +3.59%  "

$ws.Cells.Item(49, 4).Value = "'126.94"
$ws.Cells.Item(49, 5).Value = "  +1.49%  "

$ws.Cells.Item(50, 4).Value = "'2.020"
$ws.Cells.Item(50, 5).Value = "  +2.30%  "

$ws.Cells.Item(51, 4).Value = "'0.07201"
$ws.Cells.Item(51, 5).Value = "  +5.44%  "
